$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray G6 cell (and anything else outside the A:C table)
$ws.Range("G1:G10").ClearContents()

# Insert a new row at position 3 for "Preetika Shetty" (fixed spelling of
# the old "Prretika Shetty" typo), shifting the remaining rows down.
$ws.Rows("3:3").Insert()

# Rewrite the data table with the corrected/updated values.
$ws.Range("A1").Value = 52501
$ws.Range("B1").Value = "Jagannath Pidaparthy"
$ws.Range("C1").Value = 70

$ws.Range("A2").Value = 52502
$ws.Range("B2").Value = "Vishal Patil"
$ws.Range("C2").Value = 82

$ws.Range("A3").Value = 52503
$ws.Range("B3").Value = "Preetika Shetty"
$ws.Range("C3").Value = 98

$ws.Range("A4").Value = 52504
$ws.Range("B4").Value = "Sagar Mishra"
$ws.Range("C4").Value = 96

$ws.Range("A5").Value = 52505
$ws.Range("B5").Value = "Shubham Mishra"
$ws.Range("C5").Value = 93

$ws.Range("A6").Value = 52506
$ws.Range("B6").Value = "Kanchan Soni"
$ws.Range("C6").Value = 95

$ws.Range("A7").Value = 52507
$ws.Range("B7").Value = "Jai Lohani"
$ws.Range("C7").Value = 92

$ws.Range("A8").Value = 52508
$ws.Range("B8").Value = "Korol Dhanda"
$ws.Range("C8").Value = 85

$ws.Range("A9").Value = 52509
$ws.Range("B9").Value = "kaustubh Srivastava"
$ws.Range("C9").Value = 97

$ws.Range("A10").Value = 52510
$ws.Range("B10").Value = "Purva Shinde"
$ws.Range("C10").Value = 80

# The row insert pushed the former row 10 down to row 11; remove that
# now-duplicated leftover row so the table ends cleanly at row 10.
$ws.Rows("11:11").Delete()

# Update the used range / selection to match the new A1:C10 table.
$ws.Range("B3").Select()
